# Auto-generated-intent script: apply scheduled-runner price/profit updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1000
$ws.Range("I18").Value = 1000
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 1000
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -716

$ws.Range("H19").Value = 1825
$ws.Range("I19").Value = 1768.5
$ws.Range("J19").Value = 1900.3334
$ws.Range("K19").Value = 1768.5
$ws.Range("L19").Value = 1900.3334
$ws.Range("M19").Value = -1593.5
$ws.Range("N19").Value = -2250.3334

$ws.Range("H51").Value = 11143.179
$ws.Range("I51").Value = 10892.333
$ws.Range("J51").Value = 11173.28
$ws.Range("K51").Value = 10892.333
$ws.Range("L51").Value = 11173.28
$ws.Range("M51").Value = -10408.333
$ws.Range("N51").Value = -12141.28

$ws.Range("H95").Value = 53720
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 53720
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 53720
$ws.Range("N95").Value = -59212

$ws.Range("H106").Value = 3313.8333
$ws.Range("I106").Value = 2332
$ws.Range("J106").Value = 4688.4
$ws.Range("K106").Value = 2332
$ws.Range("L106").Value = 4688.4
$ws.Range("M106").Value = -1701
$ws.Range("N106").Value = -5950.4

$ws.Range("H107").Value = 593.1
$ws.Range("I107").Value = 632.44446
$ws.Range("J107").Value = 239
$ws.Range("K107").Value = 632.44446
$ws.Range("L107").Value = 239
$ws.Range("M107").Value = 1287.55554
$ws.Range("N107").Value = -4079

$ws.Range("H112").Value = 42204.766
$ws.Range("I112").Value = 78585.69500000001
$ws.Range("J112").Value = 29758.658
$ws.Range("K112").Value = 235757.085
$ws.Range("L112").Value = 89275.974
$ws.Range("M112").Value = -234649.085
$ws.Range("N112").Value = -91491.974

$ws.Range("H113").Value = 2999
$ws.Range("I113").Value = 2999
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2999
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 255

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").ClearContents()
$ws.Range("N24").Value = 0

$ws.Range("H32").Value = 3301.975
$ws.Range("I32").Value = 2345.257
$ws.Range("J32").Value = 9999
$ws.Range("K32").Value = 2345.257
$ws.Range("L32").Value = 9999
$ws.Range("M32").Value = -2058.257
$ws.Range("N32").Value = -10573

$ws.Range("H45").Value = 530242
$ws.Range("I45").Value = 591741.1
$ws.Range("J45").Value = 7499.5
$ws.Range("K45").Value = 591741.1
$ws.Range("L45").Value = 7499.5
$ws.Range("M45").Value = -591364.1
$ws.Range("N45").Value = -8253.5

$ws.Range("H96").Value = 63447.668
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 63447.668
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 63447.668
$ws.Range("N96").Value = -68939.66800000001

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").ClearContents()
$ws.Range("N100").Value = 0

$ws.Range("H102").Value = 1744.375
$ws.Range("I102").Value = 1996
$ws.Range("J102").Value = 989.5
$ws.Range("K102").Value = 1996
$ws.Range("L102").Value = 989.5
$ws.Range("M102").Value = -374
$ws.Range("N102").Value = -4233.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2708.682
$ws.Range("I86").Value = 2765.1333
$ws.Range("J86").Value = 2587.7144
$ws.Range("K86").Value = 2765.1333
$ws.Range("L86").Value = 2587.7144
$ws.Range("M86").Value = -1642.1333
$ws.Range("N86").Value = -4833.7144

$ws.Range("H89").Value = 2708.682
$ws.Range("I89").Value = 2765.1333
$ws.Range("J89").Value = 2587.7144
$ws.Range("K89").Value = 13825.6665
$ws.Range("L89").Value = 12938.572
$ws.Range("M89").Value = -8209.666499999999
$ws.Range("N89").Value = -24170.572

$ws.Range("H105").Value = 2591
$ws.Range("I105").Value = 2531
$ws.Range("J105").Value = 3011
$ws.Range("K105").Value = 2531
$ws.Range("L105").Value = 3011
$ws.Range("M105").Value = -784
$ws.Range("N105").Value = -6505

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 3485.2727
$ws.Range("I5").Value = 1322
$ws.Range("J5").Value = 5288
$ws.Range("K5").Value = 1322
$ws.Range("L5").Value = 5288
$ws.Range("M5").Value = -1210
$ws.Range("N5").Value = -5512

$ws.Range("H6").Value = 277
$ws.Range("I6").Value = 213.42857
$ws.Range("J6").Value = 499.5
$ws.Range("K6").Value = 213.42857
$ws.Range("L6").Value = 499.5
$ws.Range("M6").Value = -100.42857
$ws.Range("N6").Value = -725.5

$ws.Range("H7").Value = 214
$ws.Range("I7").Value = 80.57143000000001
$ws.Range("J7").Value = 347.42856
$ws.Range("K7").Value = 80.57143000000001
$ws.Range("L7").Value = 347.42856
$ws.Range("M7").Value = 32.42856999999999
$ws.Range("N7").Value = -573.4285600000001

$ws.Range("H13").Value = 26769.666
$ws.Range("I13").Value = 300
$ws.Range("J13").Value = 40004.5
$ws.Range("K13").Value = 300
$ws.Range("L13").Value = 40004.5
$ws.Range("M13").Value = -161
$ws.Range("N13").Value = -40282.5

$ws.Range("H15").Value = 1836
$ws.Range("I15").Value = 499.5
$ws.Range("J15").Value = 4509
$ws.Range("K15").Value = 499.5
$ws.Range("L15").Value = 4509
$ws.Range("M15").Value = -329.5
$ws.Range("N15").Value = -4849

$ws.Range("H22").Value = 20293.2
$ws.Range("I22").Value = 33561
$ws.Range("J22").Value = 391.5
$ws.Range("K22").Value = 33561
$ws.Range("L22").Value = 391.5
$ws.Range("M22").Value = -33211
$ws.Range("N22").Value = -1091.5

$ws.Range("H86").Value = 3399.4
$ws.Range("I86").Value = 3399.4
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3399.4
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -2276.4

$ws.Range("H89").Value = 3399.4
$ws.Range("I89").Value = 3399.4
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 16997
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -11381

$ws.Range("H94").Value = 2924.25
$ws.Range("I94").Value = 2699
$ws.Range("J94").Value = 2999.3333
$ws.Range("K94").Value = 2699
$ws.Range("L94").Value = 2999.3333
$ws.Range("M94").Value = -2248
$ws.Range("N94").Value = -3901.3333

$ws.Range("H134").Value = 10419729
$ws.Range("I134").Value = 12502979
$ws.Range("J134").Value = 3478
$ws.Range("K134").Value = 37508937
$ws.Range("L134").Value = 10434
$ws.Range("M134").Value = -37506402
$ws.Range("N134").Value = -15504

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 271.30768
$ws.Range("I12").Value = 246.71428
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 740.14284
$ws.Range("L12").Value = 900
$ws.Range("M12").Value = -567.14284
$ws.Range("N12").Value = -1246

$ws.Range("H14").Value = 255.5
$ws.Range("I14").Value = 255.5
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 766.5
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -593.5

$ws.Range("H70").Value = 8819.866
$ws.Range("I70").Value = 4754.364
$ws.Range("J70").Value = 20000
$ws.Range("K70").Value = 14263.092
$ws.Range("L70").Value = 60000
$ws.Range("M70").Value = -13948.092
$ws.Range("N70").Value = -60630

$ws.Range("H73").Value = 8819.866
$ws.Range("I73").Value = 4754.364
$ws.Range("J73").Value = 20000
$ws.Range("K73").Value = 14263.092
$ws.Range("L73").Value = 60000
$ws.Range("M73").Value = -13171.092
$ws.Range("N73").Value = -62184

$ws.Range("H75").Value = 3240
$ws.Range("I75").Value = 3240
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 9720
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -8722

$ws.Range("H78").Value = 3240
$ws.Range("I78").Value = 3240
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 29160
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -24168

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3097.0527
$ws.Range("I80").Value = 2614.2
$ws.Range("J80").Value = 3633.5557
$ws.Range("K80").Value = 2614.2
$ws.Range("L80").Value = 3633.5557
$ws.Range("M80").Value = -1616.2
$ws.Range("N80").Value = -5629.5557

$ws.Range("H83").Value = 3097.0527
$ws.Range("I83").Value = 2614.2
$ws.Range("J83").Value = 3633.5557
$ws.Range("K83").Value = 13071
$ws.Range("L83").Value = 18167.7785
$ws.Range("M83").Value = -8079
$ws.Range("N83").Value = -28151.7785

$ws.Range("H122").Value = 40859.062
$ws.Range("I122").Value = 51833.043
$ws.Range("J122").Value = 7937.125
$ws.Range("K122").Value = 155499.129
$ws.Range("L122").Value = 23811.375
$ws.Range("M122").Value = -153049.129
$ws.Range("N122").Value = -28711.375

$ws.Range("H126").Value = 7578.5
$ws.Range("I126").Value = 9063
$ws.Range("J126").Value = 3718.8
$ws.Range("K126").Value = 27189
$ws.Range("L126").Value = 11156.4
$ws.Range("M126").Value = -24719
$ws.Range("N126").Value = -16096.4

$ws.Range("H132").Value = 5437101
$ws.Range("I132").Value = 5437101
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 16311303
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -16308773

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 3000
$ws.Range("N68").Value = -4498

$ws.Range("H71").Value = 3000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 15000
$ws.Range("N71").Value = -22488

$ws.Range("H82").Value = 1666
$ws.Range("I82").Value = 1759.5
$ws.Range("J82").Value = 1198.5
$ws.Range("K82").Value = 1759.5
$ws.Range("L82").Value = 1198.5
$ws.Range("M82").Value = -1398.5
$ws.Range("N82").Value = -1920.5

$ws.Range("H85").Value = 1666
$ws.Range("I85").Value = 1759.5
$ws.Range("J85").Value = 1198.5
$ws.Range("K85").Value = 1759.5
$ws.Range("L85").Value = 1198.5
$ws.Range("M85").Value = -511.5
$ws.Range("N85").Value = -3694.5

$ws.Range("H132").Value = 53348236
$ws.Range("I132").Value = 60016212
$ws.Range("J132").Value = 4400
$ws.Range("K132").Value = 180048636
$ws.Range("L132").Value = 13200
$ws.Range("M132").Value = -180046106
$ws.Range("N132").Value = -18260

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1349.6428
$ws.Range("I81").Value = 1407.9166
$ws.Range("J81").Value = 1000
$ws.Range("K81").Value = 2815.8332
$ws.Range("L81").Value = 2000
$ws.Range("M81").Value = -1754.8332
$ws.Range("N81").Value = -4122

$ws.Range("H84").Value = 1349.6428
$ws.Range("I84").Value = 1407.9166
$ws.Range("J84").Value = 1000
$ws.Range("K84").Value = 14079.166
$ws.Range("L84").Value = 10000
$ws.Range("M84").Value = -8775.166000000001
$ws.Range("N84").Value = -20608

$ws.Range("H95").Value = 4672
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 4672
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 4672
$ws.Range("N95").Value = -10164

$ws.Range("H105").Value = 5000
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 5000
$ws.Range("N105").Value = -11988
